# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 6
    4  = 0
    5  = 2
    6  = 8
    7  = 5
    8  = 5
    9  = 6
    10 = 9
    11 = 2
    12 = 5
    13 = 3
    14 = 6
    15 = 5
    16 = 2
    17 = 9
    18 = 12
    19 = 5
    20 = 1
    21 = 8
    22 = 3
    23 = 3
    24 = 6
    25 = 4
    26 = 4
    27 = 4
    28 = 1
    29 = 3
    30 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
